$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old placeholder/"total" row 24 (formulas referencing an
# all-zero entry) entirely so the row becomes empty and drops out of the
# sheet's XML.
$ws.Range("E24").Clear()
$ws.Range("G24").Clear()
$ws.Range("H24").Clear()

# Populate what used to be a blank spacer row (25) with a new data entry.
$ws.Range("A25").Value = 45281.946527777778
$ws.Range("A25").NumberFormat = "m/d/yy h:mm"

$ws.Range("B25").Value = 1780736
$ws.Range("C25").Value = 600064
$ws.Range("D25").Value = 392192

$ws.Range("E25").Formula = "=SUM(B25:D25)"

$ws.Range("F25").Value = 146996204

$ws.Range("G25").Formula = "=1-(E25/E2)"
$ws.Range("G25").NumberFormat = "0.00%"

$ws.Range("H25").Formula = "=1-(F25/F2)"
$ws.Range("H25").NumberFormat = "0.00%"

$ws.Range("I25").Value = "v0.0.4, prerelease 1"

# Update the saved selection to match the author's final cursor position.
$ws.Range("J25").Select()
